$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert two new rows at position 4, pushing the rest of the metadata down.
$ws.Rows("4:5").Insert()

# New "dataset.preview.table" entry.
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

# New "dataset.preview.line" entry.
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# Formatting to match the other metadata rows plus wrap + taller rows so the
# multi-line preview scripts are fully visible.
$ws.Rows("4:5").RowHeight = 120
$ws.Range("A4:B5").VerticalAlignment = -4108
$ws.Range("A4:B5").WrapText = $true

$ws.Range("C8").Select()
